$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 74-84: a vertical-line/marker series appended below the
# existing chart-source data (columns A-C).

$data = @(
    @{ Row = 74; A = 0; B = 6.553859863;  C = 0.08916913119583336 },
    @{ Row = 75; A = 0; B = 6.467245605;  C = -0.1043174339261903 },
    @{ Row = 76; A = 0; B = 5.688794922;  C = 0.1656417864690476 },
    @{ Row = 77; A = 0; B = 6.06337793;   C = 5.963314870000002 },
    @{ Row = 78; A = 0; B = 6.553859863 },
    @{ Row = 79; A = 0; B = 5.688794922 },
    @{ Row = 80; A = 0; B = 6.553859863 },
    @{ Row = 81; A = 0; B = 6.467245605 },
    @{ Row = 82; A = 0; B = 6.553859863 },
    @{ Row = 83; A = 0; B = 6.553859863 },
    @{ Row = 84; A = 0; B = 6.553859863 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    if ($entry.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $entry.C
    }
}

# Row 84 also carries an empty text marker in column C (leading apostrophe
# forces Excel to store it as an empty text value instead of a blank cell);
# reapply the default style so only the value type changes.
$ws.Cells.Item(84, 3).Value = "'"
$ws.Cells.Item(84, 3).Style = "Normal"
